# Insert a "Saudi Arabia" row into each of the three country blocks
# (Global / High-income / International) on the "along"-by-country sheet,
# right before the existing "USA" row in each block, and fill it with its
# own mean / CI_low / CI_high values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert blank rows from the top of the sheet down; each insertion pushes
# everything below it (including the still-to-be-processed USA rows) down
# by one, so the row numbers below account for the rows already inserted.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(39).Insert()

# Global block - Saudi Arabia
$ws.Range("A13").Value = "Saudi Arabia"
$ws.Range("B13").Value = 84.0970246035313
$ws.Range("C13").Value = 80.0478220275675
$ws.Range("D13").Value = 88.1462271794951
$ws.Range("E13").Value = "**Global**:<br>Implemented by<br>All other countries"

# High-income block - Saudi Arabia
$ws.Range("A26").Value = "Saudi Arabia"
$ws.Range("B26").Value = 82.9596895936391
$ws.Range("C26").Value = 78.8773693048814
$ws.Range("D26").Value = 87.0420098823967
$ws.Range("E26").Value = "**High-income**:<br>All other HICs and<br>not some MICs (such as China)"

# International block - Saudi Arabia
$ws.Range("A39").Value = "Saudi Arabia"
$ws.Range("B39").Value = 82.1857770240346
$ws.Range("C39").Value = 78.2210970492523
$ws.Range("D39").Value = 86.1504569988168
$ws.Range("E39").Value = "**International**:<br>Some countries (e.g. EU, UK, Brazil)<br>and not others (e.g. U.S., China)"
